$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Nothern" -> "Northern" typo in the region column (column E) for
# every province row that had the misspelling. These rows all share the
# same shared-string entry, so rewriting each cell's value updates the
# whole table consistently (Excel will drop the now-unused "Nothern"
# shared string and add a single new "Northern" entry).
$rowsToFix = @(10, 11, 20, 21, 24, 33, 42, 48, 76)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 5).Value = "Northern"
}

# Match the author's final selection in the saved workbook.
$ws.Range("D7").Select()
